$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.241.32"
$ws.Range("E2").Value = "  -3.05%  "

$ws.Range("D3").Value = "2.463.33"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'310.42"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").Value = "'93.41"
$ws.Range("E6").Value = "  -6.36%  "

$ws.Range("D7").Value = "'0.548"
$ws.Range("E7").Value = "  -2.75%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  -4.79%  "

$ws.Range("D10").Value = "'33.22"
$ws.Range("E10").Value = "  -6.60%  "

$ws.Range("D11").Value = "'0.0774"
$ws.Range("E11").Value = "  -3.44%  "

$ws.Range("D12").Value = "'0.107"
$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").Value = "'6.95"
$ws.Range("E13").Value = "  -3.94%  "

$ws.Range("D14").Value = "2.842.45"
$ws.Range("E14").Value = "  -2.20%  "

$ws.Range("D15").Value = "2.474.59"
$ws.Range("E15").Value = "  -6.17%  "

$ws.Range("D16").Value = "'14.54"
$ws.Range("E16").Value = "  -5.66%  "

$ws.Range("D17").Value = "'0.776"
$ws.Range("E17").Value = "  -3.23%  "

$ws.Range("D18").Value = "41.216.50"
$ws.Range("E18").Value = "  -3.04%  "

$ws.Range("D19").Value = "'6.29"
$ws.Range("E19").Value = "  -6.63%  "

$ws.Range("D20").Value = "0.0₃0915"
$ws.Range("E20").Value = "  -3.15%  "

$ws.Range("D21").Value = "'11.21"

$ws.Range("D22").Value = "'67.96"
$ws.Range("E22").Value = "  -1.54%  "

$ws.Range("D23").Value = "'235.23"
$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("D24").Value = "'2.75"
$ws.Range("E24").Value = "  -3.53%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "'1.90"
$ws.Range("E26").Value = "  -5.75%  "

$ws.Range("D27").Value = "'23.80"
$ws.Range("E27").Value = "  -5.15%  "

$ws.Range("D28").Value = "'2.20"
$ws.Range("E28").Value = "  -5.72%  "

$ws.Range("D29").Value = "'9.57"
$ws.Range("E29").Value = "  -4.85%  "

$ws.Range("D30").Value = "'35.49"
$ws.Range("E30").Value = "  -8.09%  "

$ws.Range("D31").Value = "'152.07"
$ws.Range("E31").Value = "  -3.30%  "

$ws.Range("D32").Value = "'5.46"
$ws.Range("E32").Value = "  -4.46%  "

$ws.Range("D33").Value = "'2.64"
$ws.Range("E33").Value = "  -6.09%  "

$ws.Range("D34").Value = "'2.56"
$ws.Range("E34").Value = "  -2.88%  "

$ws.Range("D35").Value = "'0.0733"
$ws.Range("E35").Value = "  -5.67%  "

$ws.Range("D36").Value = "'2.99"
$ws.Range("E36").Value = "  -5.21%  "

$ws.Range("D37").Value = "'16.94"
$ws.Range("E37").Value = "  -4.80%  "

$ws.Range("D38").Value = "'1.87"
$ws.Range("E38").Value = "  -5.71%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.113"
$ws.Range("E39").Value = "  -3.40%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.102"
$ws.Range("E40").Value = "  -7.56%  "

$ws.Range("D41").Value = "'4.21"
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'20.21"
$ws.Range("E43").Value = "  -6.74%  "

$ws.Range("D44").Value = "1.979.58"
$ws.Range("E44").Value = "  -0.79%  "

$ws.Range("D45").Value = "'0.0282"
$ws.Range("E45").Value = "  -5.35%  "

$ws.Range("D46").Value = "'3.01"
$ws.Range("E46").Value = "  -6.69%  "

$ws.Range("D47").Value = "'8.60"
$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("D48").Value = "'69.66"
$ws.Range("E48").Value = "  -3.19%  "

$ws.Range("D49").Value = "'96.23"
$ws.Range("E49").Value = "  -4.01%  "

$ws.Range("D50").Value = "'73.95"
$ws.Range("E50").Value = "  -5.95%  "

$ws.Range("D51").Value = "'0.176"
$ws.Range("E51").Value = "  -6.45%  "
